# Auto-generated market price data refresh for Goblin_Profits workbook.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H51" = 2483.3333
    "L51" = 2400
    "N51" = -3368
    "J51" = 2400
    "L69" = 3000000000
    "J69" = 1000000000
    "H69" = 500002020
    "N69" = -3000001748
    "L70" = 200014728
    "N70" = -200015268
    "J70" = 66671576
    "H70" = 30558528
    "N72" = -9000008736
    "J72" = 1000000000
    "H72" = 500002020
    "L72" = 9000000000
    "N73" = -200016600
    "J73" = 66671576
    "L73" = 200014728
    "H73" = 30558528
    "L80" = 2583.75
    "J80" = 861.25
    "H80" = 678.3
    "N80" = -4579.75
    "N83" = -17735.25
    "J83" = 861.25
    "L83" = 7751.25
    "H83" = 678.3
    "K92" = 934.44446
    "H92" = 934.44446
    "I92" = 934.44446
    "M92" = 313.55554
    "K100" = 2775.7778
    "I100" = 2775.7778
    "M100" = -2234.7778
    "H100" = 4559.185
    "J129" = 2666.6667
    "H129" = 1176.5
    "L129" = 8000.000100000001
    "N129" = -18000.0001
    "K135" = 12817.8
    "M135" = -10282.8
    "H135" = 1457.875
    "I135" = 1424.2
    "M137" = -35305.5
    "K137" = 37855.5
    "L137" = 13676.7
    "J137" = 4558.9
    "I137" = 12618.5
    "N137" = -18776.7
    "H137" = 9260.333000000001
    "M138" = 1518.0769
    "H138" = 2529.1462
    "I138" = 1207.3077
    "K138" = 3621.9231
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "I32" = 3297.7
    "K32" = 3297.7
    "H32" = 3095.8372
    "M32" = -3010.7
    "I45" = 1590
    "M45" = -1213
    "H45" = 1590
    "K45" = 1590
    "H61" = 3914.3333
    "I61" = 3335.7646
    "K61" = 3335.7646
    "M61" = -3123.7646
    "H102" = 2505.3225
    "M102" = -106.3334
    "K102" = 1728.3334
    "I102" = 1728.3334
    "I110" = 915.1539
    "M110" = 1129.8461
    "K110" = 915.1539
    "H110" = 1059.8
    "M132" = -2851.4
    "K132" = 5381.4
    "H132" = 2102.1853
    "I132" = 1793.8
    "M136" = -7457.293799999999
    "I136" = 3335.7646
    "K136" = 10007.2938
    "H136" = 3914.3333
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "K80" = 203.28572
    "M80" = 794.71428
    "I80" = 203.28572
    "L80" = 431
    "J80" = 431
    "H80" = 376.0345
    "N80" = -2427
    "N83" = -12139
    "J83" = 431
    "I83" = 203.28572
    "L83" = 2155
    "H83" = 376.0345
    "M83" = 3975.5714
    "K83" = 1016.4286
    "K105" = 1274.579
    "L105" = 3664.2856
    "J105" = 3664.2856
    "M105" = 472.421
    "N105" = -7158.2856
    "H105" = 1917.9615
    "I105" = 1274.579
    "M134" = -8100
    "H134" = 3448.8235
    "I134" = 3545
    "K134" = 10635
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H7" = 75.47369
    "K7" = 50.6
    "J7" = 168.75
    "I7" = 50.6
    "M7" = 62.4
    "L7" = 168.75
    "N7" = -394.75
    "M134" = -4767.706200000001
    "H134" = 2386.1428
    "I134" = 2434.2354
    "K134" = 7302.706200000001
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "I19" = 45
    "H19" = 45
    "K19" = 135
    "L19" = 0
    "J19" = 0
    "M19" = 39
    "H68" = 5496.75
    "K68" = 11980.5
    "M68" = -11169.5
    "I68" = 3993.5
    "I71" = 3993.5
    "M71" = -31885.5
    "K71" = 35941.5
    "H71" = 5496.75
    "J93" = 23527.5
    "L93" = 70582.5
    "N93" = -74326.5
    "H93" = 13754.429
    "N107" = -6596.7273
    "H107" = 2013.0625
    "L107" = 2756.7273
    "J107" = 918.9091
    "K113" = 1857.6924
    "L113" = 4923.428400000001
    "N113" = -9263.428400000001
    "H113" = 976.9
    "J113" = 1641.1428
    "I113" = 619.2308
    "M113" = 312.3075999999999
    "I121" = 552.5
    "H121" = 700.9474
    "M121" = -347.5
    "N121" = -5486.2858
    "J121" = 955.4286
    "L121" = 2866.2858
    "K121" = 1657.5
    "H140" = 2323.182
    "I140" = 1832
    "M140" = -316
    "K140" = 5496
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("N19")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "I107" = 1025
    "H107" = 1025
    "K107" = 1025
    "M107" = 895
    "M126" = -4186.400000000001
    "H126" = 2977.5557
    "J126" = 3926
    "L126" = 11778
    "N126" = -16718
    "I126" = 2218.8
    "K126" = 6656.400000000001
    "N132" = -26054
    "M132" = -7188.600199999999
    "K132" = 9718.600199999999
    "H132" = 3474.4375
    "I132" = 3239.5334
    "L132" = 20994
    "J132" = 6998
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 2912.25
    "K7" = 2912.25
    "I7" = 2912.25
    "M7" = -2800.25
    "M82" = -1986.6
    "J82" = 4153.4546
    "L82" = 4153.4546
    "H82" = 2988.3872
    "N82" = -4875.4546
    "K82" = 2347.6
    "I82" = 2347.6
    "N85" = -6649.4546
    "K85" = 2347.6
    "M85" = -1099.6
    "I85" = 2347.6
    "J85" = 4153.4546
    "H85" = 2988.3872
    "L85" = 4153.4546
    "M93" = -2088.8823
    "K93" = 3336.8823
    "I93" = 3336.8823
    "H93" = 4320.3335
    "N122" = -16375
    "L122" = 11475
    "K122" = 10387.125
    "J122" = 3825
    "H122" = 3583.25
    "M122" = -7937.125
    "I122" = 3462.375
    "M126" = -6266.75
    "H126" = 2912.25
    "I126" = 2912.25
    "K126" = 8736.75
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H81" = 1364.7142
    "L81" = 4532.6666
    "N81" = -6654.6666
    "J81" = 2266.3333
    "J84" = 2266.3333
    "H84" = 1364.7142
    "L84" = 22663.333
    "N84" = -33271.333
    "L125" = 42715
    "N125" = -52555
    "J125" = 42715
    "H125" = 42715
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
